# "Mise à jour de l'application" - append 11 new GPS training rows (2025-10-21,
# "J+3" entrainement, period "Global") for the players tracked that session.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the template row (717) down to the new rows (718:728)
$ws.Range("A717:V717").Copy() | Out-Null
$ws.Range("A718:V728").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Row 718
$ws.Range("A718").Value = 'Entrainement'
$ws.Range("B718").Value = 45951
$ws.Range("C718").Value = 'Global'
$ws.Range("D718").Value = 'J+3'
$ws.Range("E718").Value = 'Amine Taiar'
$ws.Range("F718").Value = 'center back'
$ws.Range("G718").Value = '00:54:45'
$ws.Range("H718").Value = 3.91
$ws.Range("I718").Value = 0.02
$ws.Range("J718").Value = 3.88
$ws.Range("K718").Value = 0.03
$ws.Range("L718").Value = 0
$ws.Range("M718").Value = 0
$ws.Range("N718").Value = 0
$ws.Range("O718").Value = 0
$ws.Range("P718").Value = 4.1900000000000004
$ws.Range("Q718").Value = 17.940000000000001
$ws.Range("R718").Value = 4.22
$ws.Range("S718").Value = 8
$ws.Range("T718").Value = 1
$ws.Range("U718").Value = 7
$ws.Range("V718").Value = 1

# Row 719
$ws.Range("A719").Value = 'Entrainement'
$ws.Range("B719").Value = 45951
$ws.Range("C719").Value = 'Global'
$ws.Range("D719").Value = 'J+3'
$ws.Range("E719").Value = 'Maé Clavel'
$ws.Range("F719").Value = 'left back'
$ws.Range("G719").Value = '01:28:58'
$ws.Range("H719").Value = 5.97
$ws.Range("I719").Value = 0.16
$ws.Range("J719").Value = 5.8
$ws.Range("K719").Value = 0.16
$ws.Range("L719").Value = 0.01
$ws.Range("M719").Value = 0
$ws.Range("N719").Value = 0
$ws.Range("O719").Value = 0
$ws.Range("P719").Value = 3.92
$ws.Range("Q719").Value = 21.73
$ws.Range("R719").Value = 4
$ws.Range("S719").Value = 21
$ws.Range("T719").Value = 0
$ws.Range("U719").Value = 26
$ws.Range("V719").Value = 2

# Row 720
$ws.Range("A720").Value = 'Entrainement'
$ws.Range("B720").Value = 45951
$ws.Range("C720").Value = 'Global'
$ws.Range("D720").Value = 'J+3'
$ws.Range("E720").Value = 'Karahali Souaré'
$ws.Range("F720").Value = 'right forward'
$ws.Range("G720").Value = '01:05:29'
$ws.Range("H720").Value = 4.1500000000000004
$ws.Range("I720").Value = 0.08
$ws.Range("J720").Value = 4.07
$ws.Range("K720").Value = 0.08
$ws.Range("L720").Value = 0
$ws.Range("M720").Value = 0
$ws.Range("N720").Value = 0
$ws.Range("O720").Value = 0
$ws.Range("P720").Value = 3.1
$ws.Range("Q720").Value = 20.51
$ws.Range("R720").Value = 5.56
$ws.Range("S720").Value = 19
$ws.Range("T720").Value = 3
$ws.Range("U720").Value = 14
$ws.Range("V720").Value = 5

# Row 721
$ws.Range("A721").Value = 'Entrainement'
$ws.Range("B721").Value = 45951
$ws.Range("C721").Value = 'Global'
$ws.Range("D721").Value = 'J+3'
$ws.Range("E721").Value = 'Omar Benyounes'
$ws.Range("F721").Value = 'center midfield'
$ws.Range("G721").Value = '01:30:12'
$ws.Range("H721").Value = 6.64
$ws.Range("I721").Value = 0.24
$ws.Range("J721").Value = 6.39
$ws.Range("K721").Value = 0.23
$ws.Range("L721").Value = 0.02
$ws.Range("M721").Value = 0
$ws.Range("N721").Value = 0
$ws.Range("O721").Value = 0
$ws.Range("P721").Value = 4.33
$ws.Range("Q721").Value = 21.76
$ws.Range("R721").Value = 4.8
$ws.Range("S721").Value = 29
$ws.Range("T721").Value = 3
$ws.Range("U721").Value = 28
$ws.Range("V721").Value = 8

# Row 722
$ws.Range("A722").Value = 'Entrainement'
$ws.Range("B722").Value = 45951
$ws.Range("C722").Value = 'Global'
$ws.Range("D722").Value = 'J+3'
$ws.Range("E722").Value = 'Romain Thunet'
$ws.Range("F722").Value = 'center back'
$ws.Range("G722").Value = '01:29:23'
$ws.Range("H722").Value = 5.99
$ws.Range("I722").Value = 0.14000000000000001
$ws.Range("J722").Value = 5.85
$ws.Range("K722").Value = 0.12
$ws.Range("L722").Value = 0.03
$ws.Range("M722").Value = 0
$ws.Range("N722").Value = 0
$ws.Range("O722").Value = 0
$ws.Range("P722").Value = 3.93
$ws.Range("Q722").Value = 24.08
$ws.Range("R722").Value = 3.98
$ws.Range("S722").Value = 30
$ws.Range("T722").Value = 0
$ws.Range("U722").Value = 17
$ws.Range("V722").Value = 4

# Row 723
$ws.Range("A723").Value = 'Entrainement'
$ws.Range("B723").Value = 45951
$ws.Range("C723").Value = 'Global'
$ws.Range("D723").Value = 'J+3'
$ws.Range("E723").Value = 'Emmanuel Valey'
$ws.Range("F723").Value = 'left forward'
$ws.Range("G723").Value = '01:26:26'
$ws.Range("H723").Value = 6.33
$ws.Range("I723").Value = 0.07
$ws.Range("J723").Value = 6.26
$ws.Range("K723").Value = 0.07
$ws.Range("L723").Value = 0
$ws.Range("M723").Value = 0
$ws.Range("N723").Value = 0
$ws.Range("O723").Value = 0
$ws.Range("P723").Value = 3.56
$ws.Range("Q723").Value = 21.33
$ws.Range("R723").Value = 5.14
$ws.Range("S723").Value = 27
$ws.Range("T723").Value = 4
$ws.Range("U723").Value = 23
$ws.Range("V723").Value = 6

# Row 724
$ws.Range("A724").Value = 'Entrainement'
$ws.Range("B724").Value = 45951
$ws.Range("C724").Value = 'Global'
$ws.Range("D724").Value = 'J+3'
$ws.Range("E724").Value = 'Levy Ndoutoume'
$ws.Range("F724").Value = 'left back'
$ws.Range("G724").Value = '01:26:44'
$ws.Range("H724").Value = 6.01
$ws.Range("I724").Value = 0.14000000000000001
$ws.Range("J724").Value = 5.87
$ws.Range("K724").Value = 0.14000000000000001
$ws.Range("L724").Value = 0
$ws.Range("M724").Value = 0
$ws.Range("N724").Value = 0
$ws.Range("O724").Value = 0
$ws.Range("P724").Value = 4.0599999999999996
$ws.Range("Q724").Value = 22.06
$ws.Range("R724").Value = 4.71
$ws.Range("S724").Value = 33
$ws.Range("T724").Value = 3
$ws.Range("U724").Value = 24
$ws.Range("V724").Value = 2

# Row 725
$ws.Range("A725").Value = 'Entrainement'
$ws.Range("B725").Value = 45951
$ws.Range("C725").Value = 'Global'
$ws.Range("D725").Value = 'J+3'
$ws.Range("E725").Value = 'Malik Boussaid'
$ws.Range("F725").Value = 'right back'
$ws.Range("G725").Value = '01:09:22'
$ws.Range("H725").Value = 5.38
$ws.Range("I725").Value = 0.12
$ws.Range("J725").Value = 5.26
$ws.Range("K725").Value = 0.12
$ws.Range("L725").Value = 0.01
$ws.Range("M725").Value = 0
$ws.Range("N725").Value = 0
$ws.Range("O725").Value = 0
$ws.Range("P725").Value = 4.0999999999999996
$ws.Range("Q725").Value = 20.83
$ws.Range("R725").Value = 4.45
$ws.Range("S725").Value = 27
$ws.Range("T725").Value = 4
$ws.Range("U725").Value = 23
$ws.Range("V725").Value = 8

# Row 726
$ws.Range("A726").Value = 'Entrainement'
$ws.Range("B726").Value = 45951
$ws.Range("C726").Value = 'Global'
$ws.Range("D726").Value = 'J+3'
$ws.Range("E726").Value = 'Yoann Martelat'
$ws.Range("F726").Value = 'center midfield'
$ws.Range("G726").Value = '01:30:29'
$ws.Range("H726").Value = 6.21
$ws.Range("I726").Value = 0.09
$ws.Range("J726").Value = 6.13
$ws.Range("K726").Value = 0.09
$ws.Range("L726").Value = 0
$ws.Range("M726").Value = 0
$ws.Range("N726").Value = 0
$ws.Range("O726").Value = 0
$ws.Range("P726").Value = 4.1100000000000003
$ws.Range("Q726").Value = 18.91
$ws.Range("R726").Value = 3.39
$ws.Range("S726").Value = 4
$ws.Range("T726").Value = 0
$ws.Range("U726").Value = 10
$ws.Range("V726").Value = 0

# Row 727
$ws.Range("A727").Value = 'Entrainement'
$ws.Range("B727").Value = 45951
$ws.Range("C727").Value = 'Global'
$ws.Range("D727").Value = 'J+3'
$ws.Range("E727").Value = 'Mattheo Haon'
$ws.Range("F727").Value = 'right back'
$ws.Range("G727").Value = '01:30:05'
$ws.Range("H727").Value = 6.61
$ws.Range("I727").Value = 0.26
$ws.Range("J727").Value = 6.34
$ws.Range("K727").Value = 0.27
$ws.Range("L727").Value = 0.01
$ws.Range("M727").Value = 0
$ws.Range("N727").Value = 0
$ws.Range("O727").Value = 0
$ws.Range("P727").Value = 4.3099999999999996
$ws.Range("Q727").Value = 21.12
$ws.Range("R727").Value = 4.1500000000000004
$ws.Range("S727").Value = 40
$ws.Range("T727").Value = 1
$ws.Range("U727").Value = 43
$ws.Range("V727").Value = 4

# Row 728
$ws.Range("A728").Value = 'Entrainement'
$ws.Range("B728").Value = 45951
$ws.Range("C728").Value = 'Global'
$ws.Range("D728").Value = 'J+3'
$ws.Range("E728").Value = 'Ilan Ihaddadene'
$ws.Range("F728").Value = 'center midfield'
$ws.Range("G728").Value = '01:26:20'
$ws.Range("H728").Value = 6.86
$ws.Range("I728").Value = 0.17
$ws.Range("J728").Value = 6.69
$ws.Range("K728").Value = 0.16
$ws.Range("L728").Value = 0.01
$ws.Range("M728").Value = 0
$ws.Range("N728").Value = 0
$ws.Range("O728").Value = 0
$ws.Range("P728").Value = 4.66
$ws.Range("Q728").Value = 21.18
$ws.Range("R728").Value = 4.9400000000000004
$ws.Range("S728").Value = 21
$ws.Range("T728").Value = 5
$ws.Range("U728").Value = 35
$ws.Range("V728").Value = 5

# Update the active selection to reflect where the user ended up after entering the new rows
$ws.Range("D731").Select() | Out-Null
